$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.145.61'
$ws.Range('E2').Value = '  -1.60%  '

$ws.Range('D3').Value = '3.531.94'
$ws.Range('E3').Value = '  -1.72%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.20'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.02%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.98'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -4.71%  '

$ws.Range('D7').Value = '3.525.24'
$ws.Range('E7').Value = '  -1.71%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.615'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -3.02%  '

$ws.Range('E9').Value = '  +0.05%  '

$ws.Range('E10').Value = '  +5.90%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.638'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -3.65%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.65'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -5.30%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000301'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.29%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.46'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -3.34%  '

$ws.Range('D15').Value = '4.106.12'
$ws.Range('E15').Value = '  -1.50%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.28'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -4.55%  '

$ws.Range('D17').Value = '3.533.88'
$ws.Range('E17').Value = '  -1.55%  '

$ws.Range('D18').Value = '69.224.66'
$ws.Range('E18').Value = '  -1.34%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.48'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.61%  '

$ws.Range('E20').Value = '  -1.31%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '536.62'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +13.19%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.03'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.17%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '19.77'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.66%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.94'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.49%  '

$ws.Range('E25').Value = '  +0.17%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '94.11'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +6.07%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.04'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.66%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.92'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -4.99%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.05'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.63%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.77'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.31%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.35'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -5.64%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.56'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +3.41%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '64.92'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.82%  '

$ws.Range('E34').Value = '  -5.09%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '572.49'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.85%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.13'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +7.31%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '38.25'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.36%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.14%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.398'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.39%  '

$ws.Range('D40').Value = '0.0₃0762'
$ws.Range('E40').Value = '  -5.70%  '

$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.37'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -4.97%  '

$ws.Range('E42').Value = '  -7.15%  '

$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.08'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.93%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.50'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +4.41%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.96'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.79%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0442'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.03%  '

$ws.Range('D47').Value = '3.168.72'
$ws.Range('E47').Value = '  -2.14%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.22'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -4.21%  '

$ws.Range('E49').Value = '  -2.62%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.20%  '

$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '135.85'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.87%  '
